# Update the cryptos worksheet with refreshed price / volume figures,
# matching the latest GitHub Actions data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into the Price column (D) as literal text. Some of
# the price strings look like plain numbers (e.g. "587.42"), and a bare
# `.Value =` assignment lets Excel's smart-entry coerce those into numeric
# cells. Briefly forcing a Text number format for the write keeps them as
# strings (matching the workbook's original inlineStr cells), then the
# format is reset back to the sheet's default style so no stray formatting
# is left behind.
function Set-PriceText($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Rows whose Coin name / Link / Price / Volume all change (rank swap) ---
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-PriceText "D30" "495.57"
$ws.Range("E30").Value = "  -2.63%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-PriceText "D31" "7.72"
$ws.Range("E31").Value = "  +0.23%  "

# --- Rows where both Price (D) and Volume(1h) (E) change ---
$priceVolumeUpdates = @{
    2  = @("67.714.34", "  +0.83%  ")
    3  = @("2.483.52",  "  +0.15%  ")
    4  = @("0.999",     "  -0.02%  ")
    5  = @("587.42",    "  +0.35%  ")
    6  = @("175.03",    "  +1.65%  ")
    11 = @("4.95",      "  +0.53%  ")
    12 = @("0.333",     "  -0.04%  ")
    13 = @("2.935.27",  "  -0.20%  ")
    14 = @("25.27",     "  -1.18%  ")
    15 = @("67.652.64", "  +0.98%  ")
    16 = @("0.0000170", "  -0.45%  ")
    17 = @("2.457.43",  "  -0.97%  ")
    18 = @("7.41",      "  -3.00%  ")
    19 = @("10.80",     "  -1.81%  ")
    20 = @("346.54",    "  -1.16%  ")
    21 = @("4.10",      "  +1.61%  ")
    23 = @("70.69",     "  +2.44%  ")
    24 = @("4.18",      "  -1.13%  ")
    26 = @("8.81",      "  -4.15%  ")
    27 = @("2.611.17",  "  +0.23%  ")
    29 = @("0.0₃0892",  "  -1.83%  ")
    34 = @("0.999",     "  -0.05%  ")
    35 = @("164.41",    "  +0.49%  ")
    38 = @("18.27",     "  +0.72%  ")
    40 = @("1.29",      "  -3.24%  ")
    43 = @("4.77",      "  -1.17%  ")
    44 = @("2.37",      "  -0.43%  ")
    45 = @("147.80",    "  +3.03%  ")
    46 = @("3.52",      "  +0.89%  ")
    47 = @("0.510",     "  -1.11%  ")
    48 = @("0.0₆0254",  "  -3.42%  ")
    50 = @("1.56",      "  -1.11%  ")
    51 = @("0.577",     "  -1.29%  ")
}

foreach ($row in $priceVolumeUpdates.Keys) {
    $vals = $priceVolumeUpdates[$row]
    Set-PriceText "D$row" $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}

# --- Rows where only Volume(1h) (E) changes ---
$volumeOnlyUpdates = @{
    8  = "  -0.40%  "
    9  = "  +4.35%  "
    10 = "  -1.41%  "
    22 = "  +0.06%  "
    25 = "  -6.67%  "
    28 = "  -0.44%  "
    32 = "  -0.16%  "
    33 = "  -0.69%  "
    36 = "  +1.60%  "
    37 = "  -0.44%  "
    39 = "  +0.01%  "
    41 = "  +1.80%  "
    42 = "  -1.69%  "
    49 = "  -0.18%  "
}

foreach ($row in $volumeOnlyUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeOnlyUpdates[$row]
}
